# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Row -> [old value, new value] (row 6 / F6 is intentionally left untouched)
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 39
    4  = 147
    5  = 56
    7  = 1386
    8  = 600
    9  = 98
    10 = 175
    11 = 122
    12 = 176
    13 = 101
    14 = 152
    15 = 138
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
